$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Arsenal) - updated prediction data
$ws.Range("B2").Value = 4.738623103850641
$ws.Range("C2").Value = 8.661485319516407
$ws.Range("D2").Value = 0.6291353594258159
$ws.Range("E2").Value = 17
$ws.Range("F2").Value = 0.06058668371832787
$ws.Range("G2").Value = 47

# Row 8 (Brentford) - updated prediction data
$ws.Range("B8").Value = 5.050884955752212
$ws.Range("C8").Value = 7.054773082942097
$ws.Range("D8").Value = 0.4071428571428571
$ws.Range("E8").Value = 15
$ws.Range("F8").Value = 0.09525703490106394
$ws.Range("G8").Value = 24
